# Reimbursement_TestData.xlsx update
# - Rename Sheet2 -> createReimbursements, make it the active sheet
# - Fix "no" -> "No" label (affects createUnits!C3 and any row referencing it)
# - createUnits!C4:C6 flip from "Yes" to "No"
# - Populate createReimbursements with header row + one full sample row
# - Re-point selections on both sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("createUnits")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2 rename ---------------------------------------------------
$ws2.Name = "createReimbursements"

# --- createUnits (sheet1) tweaks -------------------------------------
# RunMode text fix (was lowercase "no")
$ws1.Cells.Item(3, 3).Value = "No"
# RunMode flips from "Yes" to "No" for the empty-field test rows
$ws1.Cells.Item(4, 3).Value = "No"
$ws1.Cells.Item(5, 3).Value = "No"
$ws1.Cells.Item(6, 3).Value = "No"

# Reset selection back to the top of the sheet and deselect the tab
$ws1.Range("A1").Select()

# --- createReimbursements (sheet2) content ----------------------------
# Header row
$ws2.Cells.Item(1, 1).Value = "TestCaseName"
$ws2.Cells.Item(1, 2).Value = "Test_Description"
$ws2.Cells.Item(1, 3).Value = "RunMode"
$ws2.Cells.Item(1, 4).Value = "Name"
$ws2.Cells.Item(1, 5).Value = "Description"
$ws2.Cells.Item(1, 6).Value = "Company"
$ws2.Cells.Item(1, 7).Value = "Applicable To"
$ws2.Cells.Item(1, 8).Value = "Unit"
$ws2.Cells.Item(1, 9).Value = "Approval Flow"
$ws2.Cells.Item(1, 10).Value = "Exceed"
$ws2.Cells.Item(1, 11).Value = "Project Code"
$ws2.Cells.Item(1, 12).Value = "Invoice"
$ws2.Cells.Item(1, 13).Value = "Attachments"
$ws2.Cells.Item(1, 14).Value = "Ledger"
$ws2.Cells.Item(1, 15).Value = "Limits"

# Sample data row
$ws2.Cells.Item(2, 1).Value = "Create Reimbursement Type with all Roles"
$ws2.Cells.Item(2, 2).Value = "Create Generic reimbursement type applicable to all"
$ws2.Cells.Item(2, 3).Value = "Yes"
$ws2.Cells.Item(2, 4).Value = "Reimb1-1"
$ws2.Cells.Item(2, 5).Value = "Applicable to All"
$ws2.Cells.Item(2, 8).Value = "Rupee"
$ws2.Cells.Item(2, 10).Value = "Yes"
$ws2.Cells.Item(2, 11).Value = "Yes"
$ws2.Cells.Item(2, 12).Value = "Yes"
$ws2.Cells.Item(2, 13).Value = "Yes"
$ws2.Cells.Item(2, 14).Value = "NO"

# Activate createReimbursements and set its selection/scroll position
$ws2.Activate()
$win = $excel.Windows.Item(1)
$win.ScrollColumn = 2
$win.ScrollRow = 1
$win.TabRatio = 0.5
$ws2.Range("H2").Select()
